# EUC_Perth_Assets.xlsx - Build Room fixes
# "Numpad not working, no pop-ups" / "27 needs a bit of work."
# Fix numpad functionality. Fix SAN text input pop-up.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # 4.2_Items
$ws2 = $wb.Worksheets.Item(2)   # 4.2_Timestamps
$ws3 = $wb.Worksheets.Item(3)   # BR_Items
$ws4 = $wb.Worksheets.Item(4)   # BR_Timestamps

# ---------------------------------------------------------------
# 4.2_Items: numpad-entry counts had been mistyped - fix the counts
# for "Desktop Mini" (row 2) and "Dock Thunderbolt G4" (row 3).
# ---------------------------------------------------------------
$ws1.Range("B2").Value = 152
$ws1.Range("C2").Value = 1152

$ws1.Range("B3").Value = 408
$ws1.Range("C3").Value = 1408

# ---------------------------------------------------------------
# 4.2_Timestamps: the numpad/pop-up bug produced two bogus log rows
# at the bottom (40 & 41). Replace them with the single real entry
# and drop the stray extra row.
# ---------------------------------------------------------------
$ws2.Range("A40").Value = "2023-12-08 00:41:11"
$ws2.Range("B40").Value = "Laptop 840 G9"
$ws2.Range("C40").Value = "Add 2"
$ws2.Range("D40").Value = ""

$ws2.Rows("41:41").Delete()

# ---------------------------------------------------------------
# BR_Items: initialise the Build Room items sheet with the same
# header row used on 4.2_Items.
# ---------------------------------------------------------------
$ws3.Range("A1").Value = "Item"
$ws3.Range("B1").Value = "LastCount"
$ws3.Range("C1").Value = "NewCount"
$ws3.Rows("1:1").RowHeight = 12.75

# ---------------------------------------------------------------
# BR_Timestamps: the SAN-number pop-up bug had left test data in
# row 1 and extra bogus rows. Reset it back to a clean header row,
# matching 4.2_Timestamps, and remove the leftover test rows.
# ---------------------------------------------------------------
$ws4.Range("A1").Value = "Timestamp"
$ws4.Range("B1").Value = "Item"
$ws4.Range("C1").Value = "Action"
$ws4.Range("D1").Value = "SAN Number"

$ws4.Rows("2:3").Delete()

# ---------------------------------------------------------------
# Restore view state: select row 1 on every sheet and make
# BR_Items (sheet 3) the active tab, as it needs attention next.
# ---------------------------------------------------------------
$ws1.Activate()
$ws1.Rows("1:1").Select()

$ws2.Activate()
$ws2.Rows("1:1").Select()

$ws4.Activate()
$ws4.Rows("1:1").Select()

$ws3.Activate()
$ws3.Rows("1:1").Select()

Write-Host "Build Room asset sheet fixes applied."
